$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 2 (shifts existing rows 2-20 down to 4-22),
# making room for the new "Screws/Driver Bits" row and the relocated header-label row.
$ws.Range("A2:A3").EntireRow.Insert()

# The inserted rows pick up row 1's bold/bordered formatting by default; the target
# workbook has these two rows back to plain (unstyled) cells, so strip that back off.
$ws.Range("A2:N3").ClearFormats()

# Row 1: replace text headers with numeric column-index header (0-13), keeps existing bold/border style.
$ws.Cells.Item(1, 1).Value = 0
$ws.Cells.Item(1, 2).Value = 1
$ws.Cells.Item(1, 3).Value = 2
$ws.Cells.Item(1, 4).Value = 3
$ws.Cells.Item(1, 5).Value = 4
$ws.Cells.Item(1, 6).Value = 5
$ws.Cells.Item(1, 7).Value = 6
$ws.Cells.Item(1, 8).Value = 7
$ws.Cells.Item(1, 9).Value = 8
$ws.Cells.Item(1, 10).Value = 9
$ws.Cells.Item(1, 11).Value = 10
$ws.Cells.Item(1, 12).Value = 11
$ws.Cells.Item(1, 13).Value = 12
$ws.Cells.Item(1, 14).Value = 13

# Row 2: new group header row.
$ws.Cells.Item(2, 1).Value = "Screws"
$ws.Cells.Item(2, 2).Value = "Driver Bits"
$ws.Cells.Item(2, 3).Value = ""
$ws.Cells.Item(2, 4).Value = ""
$ws.Cells.Item(2, 5).Value = ""
$ws.Cells.Item(2, 6).Value = ""
$ws.Cells.Item(2, 7).Value = ""
$ws.Cells.Item(2, 8).Value = ""
$ws.Cells.Item(2, 9).Value = ""
$ws.Cells.Item(2, 10).Value = ""
$ws.Cells.Item(2, 11).Value = ""
$ws.Cells.Item(2, 12).Value = ""
$ws.Cells.Item(2, 13).Value = ""
$ws.Cells.Item(2, 14).Value = ""

# Row 3: the original column-label header text, now un-styled (plain row).
$ws.Cells.Item(3, 1).Value = "Lg."
$ws.Cells.Item(3, 2).Value = "Threading"
$ws.Cells.Item(3, 3).Value = "HeadDia."
$ws.Cells.Item(3, 4).Value = "Head Ht."
$ws.Cells.Item(3, 5).Value = "DriveSize"
$ws.Cells.Item(3, 6).Value = "TensileStrength, psi"
$ws.Cells.Item(3, 7).Value = "SpecificationsMet"
$ws.Cells.Item(3, 8).Value = "Pkg.Qty."
$ws.Cells.Item(3, 9).Value = ""
$ws.Cells.Item(3, 10).Value = "Pkg."
$ws.Cells.Item(3, 11).Value = ""
$ws.Cells.Item(3, 12).Value = "Each"
$ws.Cells.Item(3, 13).Value = ""
$ws.Cells.Item(3, 14).Value = ""